# Fix id error on diagram
#  1) "Software Component" sheet/table: add an "Out of scope" column indicating
#     whether each software component is in/out of scope.
#  2) "Trust Boundaries" sheet: fix a data-entry mix-up where the "Enablers"
#     and "Enablers2" boundary identifiers were swapped between rows 2 and 3.

$wb = $excel.ActiveWorkbook

# --- Sheet "Software Component": add "Out of scope" column to Table1 ---
$ws1 = $wb.Worksheets.Item("Software Component")
$lo1 = $ws1.ListObjects.Item(1)
$null = $lo1.ListColumns.Add()

$ws1.Range("D1").Value2 = "Out of scope"
$ws1.Range("D2").Value2 = "No"
$ws1.Range("D3").Value2 = "Yes"
$ws1.Range("D4").Value2 = "Yes"
$ws1.Range("D5").Value2 = "Yes"
$ws1.Range("D6").Value2 = "Yes"

$ws1.Columns.Item(4).ColumnWidth = 13.67

# --- Sheet "Trust Boundaries": swap the mixed-up Enablers / Enablers2 ids ---
$ws3 = $wb.Worksheets.Item("Trust Boundaries")
$ws3.Range("A2").Value2 = "Enablers2"
$ws3.Range("A3").Value2 = "Enablers"
